$d = $word.ActiveDocument

# --- 1) Grab the "Meta description" label run's formatting (an empty run
#        followed by the bold "Meta description" run) from the paragraph
#        right after the H1 title. We'll reuse it for the new paragraph
#        below, then delete the original paragraph afterwards (deleting it
#        first would shift character positions out from under any range
#        that still points at this paragraph).
$metaPara = $d.Paragraphs.Item(2)
$labelRange = $metaPara.Range.Duplicate
$labelRange.Collapse(1)
$labelRange.MoveEnd(1, 16)
$labelFormattedText = $labelRange.FormattedText

# --- 2) Append a new bold paragraph (re-using the label run's formatting)
#        right before the closing "Prompt for DALLE" paragraph, with the
#        site title text.
$n = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($n - 1)
$secondToLast.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($n)
$newPara.Range.Style = -1
$newPara.Range.FormattedText = $labelFormattedText
$newPara.Range.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "Play Forgotten Fable for Free - Stunning Graphics & Exciting Features", 2)

# --- 3) Swap the closing DALLE image-prompt paragraph's text for the old
#        meta description copy (keeping its italic formatting intact).
$d.Content.Find.Execute("Prompt for DALLE: Create a colorful cartoon-style feature image that captures the essence of the game " + [char]34 + "Forgotten Fable" + [char]34 + ". The image should prominently feature a happy Maya warrior with glasses. The warrior should be shown alongside the game's villain characters, including Baba Yaga, Koschei the Immortal, Vodyanoy, Leshy, and Gorynych. The background should showcase the game's fantasy theme with dragons, magic chests, and runes. The image should be captivating and draw the attention of viewers, with attention paid to the smallest details to showcase the Evoplay's excellent graphics.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the Forgotten Fable slot game with its unique fantasy theme, stunning graphics, and exciting special features. Play for free and win big!", 2)

# --- 4) Finally, remove the whole "Meta description: ..." paragraph that
#        sits right after the title heading (it's still paragraph #2; none
#        of the above touched that part of the document).
$metaPara2 = $d.Paragraphs.Item(2)
$metaPara2.Range.Delete()

Write-Output "done"
